# Update the ticket-price column (G, "最低票价") on both the "展览" and
# "全部类型" sheets. These used to store a raw numeric "cents" price (or 0
# when unknown); they now hold human-readable text - either the price in
# whole yuan, or a status string for listings that aren't purchasable.
#
# A leading apostrophe forces Excel to store a numeric-looking value (e.g.
# "65") as text instead of silently re-parsing it back into a number.
$gValues = @{
    2  = "'不可售"
    3  = "'65"
    4  = "'已售罄"
    5  = "'258"
    6  = "'55"
    7  = "'55"
    8  = "'55"
    9  = "'65"
    10 = "'168"
    11 = "'60"
    12 = "'70"
    13 = "'39"
    14 = "'60"
    15 = "'65"
    16 = "'68"
    17 = "'60"
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $gValues.Keys) {
        $ws.Cells.Item($row, 7).Value = $gValues[$row]
    }

    # Row 15's "想去人数" (F) count also moved from 1104 to 1105.
    $ws.Cells.Item(15, 6).Value = 1105
}
